$d = $word.ActiveDocument

$newText = @"
1. Настоящий межевой план подготовлен в соответствии с документацией по планировке территории объекта «Строительство скоростной автомобильной дороги Москва-Санкт-Петербург на участке км 58 – км 684 (с последующей эксплуатацией на платной основе)». 7 этап км 543 – км 646». утвержденной Распоряжением Федерального дорожного агентства (РОСАВТОДОР) № 907-р от 13.05.2014г. Следовательно были соблюдены все этапы согласования образуемых земельных участков со смежниками.
2. Межевой план подготовлен  всоответствии с проектными документациями лесного участка, утвержденных распоряжениями Комитета по природным ресурсам Ленинградской области № 566 и 565 от 12.04.2018г.   
3. Сформированные земельные участки, представленные в данном межевом плане, имеют порядковые номера 6, 14 в проекте межевания и подлежит изъятию для нужд Российской Федерации в соответствии с Распоряжением Федерального дорожного агенста  №907-р от 13.05.2014  "Об изъятии для нужд Российской Федерации земельных участков в целях обеспечения реализации проекта «Строительство скоростной автомобильной дороги Москва-Санкт-Петербург на участке км 58 – км 684 (с последующей эксплуатацией на платной основе)». 7 этап км 543 – км 646. 
4. В соответствии с п.5 ст.36, п.8 ст.37 Федеральный закон от 17.07.2009 N 145-ФЗ Государственная компания вправе обращаться с заявлениями о проведении государственного кадастрового учета земельных участков, предназначенных для размещения автомобильных дорог, деятельность по организации строительства или реконструкции которых осуществляет Государственная компания.
5.Согласно пункту 32 Приказа № 921 от 08.12.2015 г. «Об утверждении формы и состава сведений межевого плана. требований к его подготовке» в реквизите «Сведения о геодезической основе. использованной при подготовке межевого плана» раздела «Исходные данные» указываются дата выполненного при проведении кадастровых работ обследования состояния наружного знака пункта и слова «сохранился». «не обнаружен» или «утрачен» в зависимости от состояния такого пункта. Т. к. в xml-схеме 6 версии нет возможности прописать данную информацию в соответствующем реквизите. электронный образ раздела «Исходные данные» в формате pdf приведен в приложении.
В связи с тем. что геодезические работы проводились в 2014г.. то сведения о сохранности геодезической основы приводятся на данный период времени.
6.Межевой план подготовлен Кадастровым инженером Поляковым Павлом Владимировичем, являющейся членом Ассоциации Саморегулируемой организации «Объединение профессионалов кадастровой деятельности» (сокращенное наименование – Ассоциация СРО «ОПКД») (уникальный реестровый номер кадастрового инженера в реестре членов Ассоциации СРО «ОПКД» № 1944). Сведения о СРО КИ: Ассоциация Саморегулируемой организации «Объединение профессионалов кадастровой деятельности» (сокращенное наименование – Ассоциация СРО «ОПКД») содержатся в государственном реестре СРО КИ (уникальный номер реестровой записи от 08.07.2016 г. № 003).
7. № регистрации в государственном реестре лиц, осуществляющих кадастровую деятельность: 34626.
8. Страховой номер индивидуального лицевого счета в системе обязательного пенсионного страхования Российской Федерации (СНИЛС): 133-227-044 15. 
9. Заказчиком работ является ООО Строительно-производственная компания «Зеленый город». Договор субподряда № КТ/09-13 от 27.09.2013г.
"@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Представленный межевой план подготовлен")) {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = $newText
        break
    }
}

Write-Output "done"
